$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column C (CPLE3) values for rows 4-33, and the resulting column O (AXIA6) totals for rows 4-35
$ws.Range("C4").Value = 3235.4861380697657
$ws.Range("O4").Value = 18511.406301685623
$ws.Range("C5").Value = 3660.143747464982
$ws.Range("O5").Value = 22879.286074147032
$ws.Range("C6").Value = 3582.3763563740481
$ws.Range("O6").Value = 21037.072244304934
$ws.Range("C7").Value = 4425.1607962207636
$ws.Range("O7").Value = 14257.130957783082
$ws.Range("C8").Value = 4210.439705807612
$ws.Range("O8").Value = 14572.640722376331
$ws.Range("C9").Value = 5058.9868636215051
$ws.Range("O9").Value = 14256.83359063529
$ws.Range("C10").Value = 5490.5886142238132
$ws.Range("O10").Value = 14916.449828221348
$ws.Range("C11").Value = 5449.6911283321651
$ws.Range("O11").Value = 14802.022900230375
$ws.Range("C12").Value = 5465.2057201785474
$ws.Range("O12").Value = 14927.501604703333
$ws.Range("C13").Value = 5657.9555381712644
$ws.Range("O13").Value = 8966.776435412432
$ws.Range("C14").Value = 5565.9368892212588
$ws.Range("O14").Value = 9170.2795180922367
$ws.Range("C15").Value = 5581.4196732150995
$ws.Range("O15").Value = 9901.0095168224398
$ws.Range("C16").Value = 5667.3724990695691
$ws.Range("O16").Value = 16924.708304675674
$ws.Range("C17").Value = 5714.844584718162
$ws.Range("O17").Value = 18304.907510996669
$ws.Range("C18").Value = 5646.4806655678176
$ws.Range("O18").Value = 19338.066657106887
$ws.Range("C19").Value = 5314.7528189505028
$ws.Range("O19").Value = 20929.557974011255
$ws.Range("C20").Value = 7991.7263392091545
$ws.Range("O20").Value = 28513.605545167349
$ws.Range("C21").Value = 4418.7495968982048
$ws.Range("O21").Value = 19008.491823315748
$ws.Range("C22").Value = 4302.338533913191
$ws.Range("O22").Value = 19847.81212405007
$ws.Range("C23").Value = -5318.5893187219672
$ws.Range("O23").Value = 20539.534716877293
$ws.Range("C24").Value = 20472.232012192413
$ws.Range("O24").Value = 20702.885848859736
$ws.Range("C25").Value = 1150.1898105245825
$ws.Range("O25").Value = 19166.341385098109
$ws.Range("C26").Value = 1000.2915388927777
$ws.Range("O26").Value = 18421.354309677256
$ws.Range("C27").Value = 1013.2949582724123
$ws.Range("O27").Value = 19402.446106210256
$ws.Range("C28").Value = 982.63393501217752
$ws.Range("O28").Value = 21169.265550508873
$ws.Range("C29").Value = 949.36586851326774
$ws.Range("O29").Value = 20989.256995824366
$ws.Range("C30").Value = 1227.873670728768
$ws.Range("O30").Value = 21125.568010647345
$ws.Range("C31").Value = 826.4587927119029
$ws.Range("O31").Value = 6548.0965375347387
$ws.Range("C32").Value = 236.94676705743296
$ws.Range("O32").Value = 5188.8253978749399
$ws.Range("C33").Value = -3197.1599030428097
$ws.Range("O33").Value = 4226.6679307624145
$ws.Range("O34").Value = 3261.4552517639413
$ws.Range("O35").Value = 1144.9852322624533

# Clear the stale O36 value (no longer carries a total in that column)
$ws.Range("O36").ClearContents()

# Restore the selection left by the editor (column C rows 4:33 were the edited range)
$ws.Activate()
$ws.Range("C4:C33").Select()
